$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 109 (shifts existing rows 109..175 down to 110..176)
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with the new observation
$ws.Range("A109").Value = 9
$ws.Range("B109").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C109").Value = "Metropolitana"
$ws.Range("D109").Value = 44518
$ws.Range("E109").Value = 13
$ws.Range("F109").Value = 100112043
$ws.Range("G109").Value = "Pepino ensalada"
$ws.Range("H109").Value = "Sin especificar"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 160
$ws.Range("K109").Value = 5000
$ws.Range("L109").Value = 7000
$ws.Range("M109").Value = 6000
$ws.Range("N109").Value = "$/caja 50 unidades"
$ws.Range("O109").Value = "Región de Arica y Parinacota"
$ws.Range("P109").Value = 120
$ws.Range("Q109").Value = 50
$ws.Range("R109").Value = "Hortaliza"
